$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Trim the English parenthetical suffix from the color-scheme names in column A (rows 8-13)
$ws.Range("A8").Value  = "森林晨曦 "
$ws.Range("A9").Value  = "寒地松林 "
$ws.Range("A10").Value = "復古大地 "
$ws.Range("A11").Value = "迷彩灰綠 "
$ws.Range("A12").Value = "深空月石"
$ws.Range("A13").Value = "經典商務"

# Row heights shrink now that the wrapped text is shorter (auto-fit result)
$ws.Range("A8:D13").EntireRow.AutoFit()

# Scroll/selection state changes recorded in the sheet view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("A7").Select()

Write-Host "done"
